$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5534.62067565243
$ws.Range("C2").Value = 17106.26726406902
$ws.Range("D2").Value = 1154.901341187198
